$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (shardulpakhare@gmail.com / Pass1234), shifting the row below it up.
$ws.Rows.Item(2).Delete()

# Update the active selection to A9, matching the post-edit cursor position.
$ws.Range("A9").Select()
